# Fruta / hortaliza, semanal
# Insert two new weekly records for "Vega Monumental Concepción - Cereza" just
# after the existing row 157, pushing the rest of the table (old rows 158-199)
# down to rows 160-201.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 158 (old row 158 -> 160, old row 159 -> 161, etc.)
$ws.Rows.Item(158).EntireRow.Insert()
$ws.Rows.Item(158).EntireRow.Insert()

# Seed the two new rows from the (now shifted) rows directly below them, so
# that all of the constant columns (Mercado ID, Mercado, Región, Codreg,
# Tipo, Producto ID, Producto, Categoría ID, Categoría, and any unchanged
# Variedad/Calidad/Unidad/Origen/Kg-unit values) and cell styles/number
# formats line up correctly. Only the used range (A:T) is copied so no
# stray formatting bleeds into unused columns.
$ws.Range("A160:T160").Copy($ws.Range("A158"))
$ws.Range("A161:T161").Copy($ws.Range("A159"))

# Row 158: new Early Burlat / Primera record (2023-11-28, serial 45258)
$ws.Cells.Item(158, 4).Value = 45258
$ws.Cells.Item(158, 11).Value = "Early Burlat"
$ws.Cells.Item(158, 13).Value = 100
$ws.Cells.Item(158, 14).Value = 11000
$ws.Cells.Item(158, 15).Value = 12000
$ws.Cells.Item(158, 16).Value = 11500
$ws.Cells.Item(158, 19).Value = 1150

# Row 159: new Royal Dawn / Primera record (2023-11-28, serial 45258)
$ws.Cells.Item(159, 4).Value = 45258
$ws.Cells.Item(159, 14).Value = 12000
$ws.Cells.Item(159, 15).Value = 13000
$ws.Cells.Item(159, 16).Value = 12500
$ws.Cells.Item(159, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(159, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(159, 19).Value = 1250
$ws.Cells.Item(159, 20).Value = 10
